$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 19
$ws_ALC.Range("H19").Value = 10618.4
$ws_ALC.Range("I19").Value = 411.6
$ws_ALC.Range("J19").Value = 14020.667
$ws_ALC.Range("K19").Value = 411.6
$ws_ALC.Range("L19").Value = 14020.667
$ws_ALC.Range("M19").Value = -236.6
$ws_ALC.Range("N19").Value = -14370.667

# ALC row 47
$ws_ALC.Range("H47").Value = 15000
$ws_ALC.Range("I47").Value = 0
$ws_ALC.Range("K47").Value = 0
$ws_ALC.Range("M47").Value = $null

# ALC row 98
$ws_ALC.Range("H98").Value = 772.1
$ws_ALC.Range("I98").Value = 539.75
$ws_ALC.Range("J98").Value = 1701.5
$ws_ALC.Range("K98").Value = 539.75
$ws_ALC.Range("L98").Value = 1701.5
$ws_ALC.Range("M98").Value = 958.25
$ws_ALC.Range("N98").Value = -4697.5

# ALC row 112
$ws_ALC.Range("H112").Value = 10379.926
$ws_ALC.Range("J112").Value = 10379.926
$ws_ALC.Range("L112").Value = 31139.778
$ws_ALC.Range("N112").Value = -33355.778

# ALC row 122
$ws_ALC.Range("H122").Value = 772.1
$ws_ALC.Range("I122").Value = 539.75
$ws_ALC.Range("J122").Value = 1701.5
$ws_ALC.Range("K122").Value = 1619.25
$ws_ALC.Range("L122").Value = 5104.5
$ws_ALC.Range("M122").Value = 830.75
$ws_ALC.Range("N122").Value = -10004.5

# ALC row 132
$ws_ALC.Range("H132").Value = 2700.8135
$ws_ALC.Range("I132").Value = 2549.9812
$ws_ALC.Range("K132").Value = 7649.943600000001
$ws_ALC.Range("M132").Value = -5119.943600000001

# ALC row 137
$ws_ALC.Range("H137").Value = 1181.6836
$ws_ALC.Range("I137").Value = 1093.8636
$ws_ALC.Range("J137").Value = 1292.0857
$ws_ALC.Range("K137").Value = 3281.5908
$ws_ALC.Range("L137").Value = 3876.2571
$ws_ALC.Range("M137").Value = -731.5907999999999
$ws_ALC.Range("N137").Value = -8976.257100000001

# ALC row 138
$ws_ALC.Range("H138").Value = 1049.1
$ws_ALC.Range("I138").Value = 563.5472
$ws_ALC.Range("J138").Value = 1596.6383
$ws_ALC.Range("K138").Value = 1690.6416
$ws_ALC.Range("L138").Value = 4789.9149
$ws_ALC.Range("M138").Value = 3449.3584
$ws_ALC.Range("N138").Value = -15069.9149

# ARM row 32
$ws_ARM.Range("H32").Value = 638025.25
$ws_ARM.Range("I32").Value = 749481
$ws_ARM.Range("J32").Value = 13873.134
$ws_ARM.Range("K32").Value = 749481
$ws_ARM.Range("L32").Value = 13873.134
$ws_ARM.Range("M32").Value = -749194
$ws_ARM.Range("N32").Value = -14447.134

# ARM row 122
$ws_ARM.Range("H122").Value = 54581.79
$ws_ARM.Range("I122").Value = 78472.30499999999
$ws_ARM.Range("J122").Value = 2819
$ws_ARM.Range("K122").Value = 235416.915
$ws_ARM.Range("L122").Value = 8457
$ws_ARM.Range("M122").Value = -232966.915
$ws_ARM.Range("N122").Value = -13357

# ARM row 132
$ws_ARM.Range("H132").Value = 2825.283
$ws_ARM.Range("I132").Value = 2963.2666
$ws_ARM.Range("J132").Value = 2645.3044
$ws_ARM.Range("K132").Value = 8889.799800000001
$ws_ARM.Range("L132").Value = 7935.9132
$ws_ARM.Range("M132").Value = -6359.799800000001
$ws_ARM.Range("N132").Value = -12995.9132

# BSM row 27
$ws_BSM.Range("I27").Value = 35000
$ws_BSM.Range("J27").Value = 69990
$ws_BSM.Range("K27").Value = 35000
$ws_BSM.Range("L27").Value = 69990
$ws_BSM.Range("N27").Value = -70374
$ws_BSM.Range("M27").Value = -34808

# BSM row 99
$ws_BSM.Range("H99").Value = 1498.7142
$ws_BSM.Range("J99").Value = 1997
$ws_BSM.Range("L99").Value = 1997
$ws_BSM.Range("N99").Value = -4993

# BSM row 105
$ws_BSM.Range("H105").Value = 17859336
$ws_BSM.Range("I105").Value = 31251724
$ws_BSM.Range("J105").Value = 2816.6667
$ws_BSM.Range("K105").Value = 31251724
$ws_BSM.Range("L105").Value = 2816.6667
$ws_BSM.Range("M105").Value = -31249977
$ws_BSM.Range("N105").Value = -6310.6667

# BSM row 134
$ws_BSM.Range("H134").Value = 1948.6046
$ws_BSM.Range("I134").Value = 1858.1578
$ws_BSM.Range("J134").Value = 2636
$ws_BSM.Range("K134").Value = 5574.4734
$ws_BSM.Range("L134").Value = 7908
$ws_BSM.Range("M134").Value = -3039.4734
$ws_BSM.Range("N134").Value = -12978

# CRP row 14
$ws_CRP.Range("H14").Value = 0
$ws_CRP.Range("I14").Value = 0
$ws_CRP.Range("K14").Value = 0
$ws_CRP.Range("M14").Value = $null

# CRP row 23
$ws_CRP.Range("H23").Value = 10347.934
$ws_CRP.Range("I23").Value = 4516.9
$ws_CRP.Range("J23").Value = 22010
$ws_CRP.Range("K23").Value = 4516.9
$ws_CRP.Range("L23").Value = 22010
$ws_CRP.Range("M23").Value = -4276.9
$ws_CRP.Range("N23").Value = -22490

# CRP row 27
$ws_CRP.Range("H27").Value = 10347.934
$ws_CRP.Range("I27").Value = 4516.9
$ws_CRP.Range("J27").Value = 22010
$ws_CRP.Range("K27").Value = 4516.9
$ws_CRP.Range("L27").Value = 22010
$ws_CRP.Range("M27").Value = -4324.9
$ws_CRP.Range("N27").Value = -22394

# CRP row 94
$ws_CRP.Range("H94").Value = 1130.8572
$ws_CRP.Range("J94").Value = 1192.6666
$ws_CRP.Range("L94").Value = 1192.6666
$ws_CRP.Range("N94").Value = -2094.6666

# CRP row 122
$ws_CRP.Range("H122").Value = 1925.4
$ws_CRP.Range("I122").Value = 1809.2858
$ws_CRP.Range("J122").Value = 1970.5555
$ws_CRP.Range("K122").Value = 5427.857400000001
$ws_CRP.Range("L122").Value = 5911.666499999999
$ws_CRP.Range("M122").Value = -2977.857400000001
$ws_CRP.Range("N122").Value = -10811.6665

# CRP row 134
$ws_CRP.Range("H134").Value = 2709.0923
$ws_CRP.Range("I134").Value = 2573.25
$ws_CRP.Range("J134").Value = 3554.3333
$ws_CRP.Range("K134").Value = 7719.75
$ws_CRP.Range("L134").Value = 10662.9999
$ws_CRP.Range("M134").Value = -5184.75
$ws_CRP.Range("N134").Value = -15732.9999

# CUL row 54
$ws_CUL.Range("H54").Value = 5407.857
$ws_CUL.Range("I54").Value = 0
$ws_CUL.Range("J54").Value = 5407.857
$ws_CUL.Range("K54").Value = 0
$ws_CUL.Range("L54").Value = 16223.571
$ws_CUL.Range("M54").Value = $null
$ws_CUL.Range("N54").Value = -17341.571

# CUL row 86
$ws_CUL.Range("H86").Value = 416.6
$ws_CUL.Range("I86").Value = 390
$ws_CUL.Range("J86").Value = 434.33334
$ws_CUL.Range("K86").Value = 1170
$ws_CUL.Range("L86").Value = 1303.00002
$ws_CUL.Range("M86").Value = 16
$ws_CUL.Range("N86").Value = -3675.00002

# CUL row 89
$ws_CUL.Range("H89").Value = 416.6
$ws_CUL.Range("I89").Value = 390
$ws_CUL.Range("J89").Value = 434.33334
$ws_CUL.Range("K89").Value = 3510
$ws_CUL.Range("L89").Value = 3909.00006
$ws_CUL.Range("M89").Value = 2418
$ws_CUL.Range("N89").Value = -15765.00006

# CUL row 92
$ws_CUL.Range("H92").Value = 700
$ws_CUL.Range("I92").Value = 900
$ws_CUL.Range("K92").Value = 2700
$ws_CUL.Range("M92").Value = -1452

# CUL row 98
$ws_CUL.Range("H98").Value = 547.8182
$ws_CUL.Range("I98").Value = 447.2
$ws_CUL.Range("J98").Value = 631.6667
$ws_CUL.Range("K98").Value = 1341.6
$ws_CUL.Range("L98").Value = 1895.0001
$ws_CUL.Range("M98").Value = 156.4000000000001
$ws_CUL.Range("N98").Value = -4891.0001

# CUL row 131
$ws_CUL.Range("H131").Value = 3077.0908
$ws_CUL.Range("I131").Value = 414.14285
$ws_CUL.Range("J131").Value = 3986.3901
$ws_CUL.Range("K131").Value = 1242.42855
$ws_CUL.Range("L131").Value = 11959.1703
$ws_CUL.Range("M131").Value = 3797.57145
$ws_CUL.Range("N131").Value = -22039.1703

# GSM row 39
$ws_GSM.Range("H39").Value = 30000
$ws_GSM.Range("J39").Value = 30000
$ws_GSM.Range("L39").Value = 30000
$ws_GSM.Range("N39").Value = -31064

# GSM row 70
$ws_GSM.Range("H70").Value = 9127.92
$ws_GSM.Range("I70").Value = 9879.9
$ws_GSM.Range("J70").Value = 6120
$ws_GSM.Range("K70").Value = 9879.9
$ws_GSM.Range("L70").Value = 6120
$ws_GSM.Range("M70").Value = -9609.9
$ws_GSM.Range("N70").Value = -6660

# GSM row 73
$ws_GSM.Range("H73").Value = 9127.92
$ws_GSM.Range("I73").Value = 9879.9
$ws_GSM.Range("J73").Value = 6120
$ws_GSM.Range("K73").Value = 9879.9
$ws_GSM.Range("L73").Value = 6120
$ws_GSM.Range("M73").Value = -8943.9
$ws_GSM.Range("N73").Value = -7992

# GSM row 126
$ws_GSM.Range("H126").Value = 1056.25
$ws_GSM.Range("I126").Value = 859
$ws_GSM.Range("J126").Value = 1253.5
$ws_GSM.Range("K126").Value = 2577
$ws_GSM.Range("L126").Value = 3760.5
$ws_GSM.Range("M126").Value = -107
$ws_GSM.Range("N126").Value = -8700.5

# GSM row 132
$ws_GSM.Range("H132").Value = 2253.4626
$ws_GSM.Range("I132").Value = 1824.6123
$ws_GSM.Range("J132").Value = 3420.889
$ws_GSM.Range("K132").Value = 5473.8369
$ws_GSM.Range("L132").Value = 10262.667
$ws_GSM.Range("M132").Value = -2943.8369
$ws_GSM.Range("N132").Value = -15322.667

# LTW row 16
$ws_LTW.Range("H16").Value = 1115.5454
$ws_LTW.Range("I16").Value = 1184.375
$ws_LTW.Range("J16").Value = 932
$ws_LTW.Range("K16").Value = 1184.375
$ws_LTW.Range("L16").Value = 932
$ws_LTW.Range("M16").Value = -1014.375
$ws_LTW.Range("N16").Value = -1272

# LTW row 68
$ws_LTW.Range("H68").Value = 1478.4717
$ws_LTW.Range("I68").Value = 1447.9778
$ws_LTW.Range("J68").Value = 1650
$ws_LTW.Range("K68").Value = 1447.9778
$ws_LTW.Range("L68").Value = 1650
$ws_LTW.Range("M68").Value = -698.9777999999999
$ws_LTW.Range("N68").Value = -3148

# LTW row 71
$ws_LTW.Range("H71").Value = 1478.4717
$ws_LTW.Range("I71").Value = 1447.9778
$ws_LTW.Range("J71").Value = 1650
$ws_LTW.Range("K71").Value = 7239.888999999999
$ws_LTW.Range("L71").Value = 8250
$ws_LTW.Range("M71").Value = -3495.888999999999
$ws_LTW.Range("N71").Value = -15738

# LTW row 136
$ws_LTW.Range("H136").Value = 3402640.2
$ws_LTW.Range("I136").Value = 1038.1389
$ws_LTW.Range("J136").Value = 12822461
$ws_LTW.Range("K136").Value = 3114.4167
$ws_LTW.Range("L136").Value = 38467383
$ws_LTW.Range("M136").Value = -564.4166999999998
$ws_LTW.Range("N136").Value = -38472483

# WVR row 107
$ws_WVR.Range("H107").Value = 950.3333
$ws_WVR.Range("I107").Value = 919.61536
$ws_WVR.Range("J107").Value = 1150
$ws_WVR.Range("K107").Value = 2758.84608
$ws_WVR.Range("L107").Value = 3450
$ws_WVR.Range("M107").Value = -838.8460800000003
$ws_WVR.Range("N107").Value = -7290

# WVR row 126
$ws_WVR.Range("H126").Value = 1180.2106
$ws_WVR.Range("I126").Value = 1200.875
$ws_WVR.Range("J126").Value = 1070
$ws_WVR.Range("K126").Value = 3602.625
$ws_WVR.Range("L126").Value = 3210
$ws_WVR.Range("M126").Value = -1132.625
$ws_WVR.Range("N126").Value = -8150

# WVR row 136
$ws_WVR.Range("H136").Value = 2039.6105
$ws_WVR.Range("I136").Value = 1974.8806
$ws_WVR.Range("J136").Value = 2194.5
$ws_WVR.Range("K136").Value = 5924.641799999999
$ws_WVR.Range("L136").Value = 6583.5
$ws_WVR.Range("M136").Value = -3374.641799999999
$ws_WVR.Range("N136").Value = -11683.5
